# #5: property boat&car done
#
# The "汽車" (car) sheet (3rd worksheet, sheet3.xml) is reshaped from a
# 7-column (A:G) table whose row 1 was actually a stray duplicate of the
# data row (no real header) into a full 14-column (A:N) table with a
# proper header row - matching the column layout already used on the
# "股票" (stock) sheet: name / capacity / owner / register_date /
# register_reason / acquire_value / property_category / category / date /
# legislator_name / legislator_id / source_file / index.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# J holds a literal "yyyy-mm-dd" string (like the other sheets' "date"
# column) - force Text format first so Excel doesn't reinterpret it as a
# date serial.
$ws.Range("J1:J2").NumberFormat = "@"

# --- Row 1: header labels (style already bold/bordered from before) -----
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Row 2: data ----------------------------------------------------------
$ws.Range("B2").Value = "NISSAN"
$ws.Range("E2").Value = "91年08月27日"
$ws.Range("G2").Value = "(超過五年）"
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Value = "2013-12-03"
$ws.Range("K2").Value = "邱議瑩"
$ws.Range("L2").Value = 913
$ws.Range("M2").Value = "tmp40191"
$ws.Range("N2").Value = 30

# --- Formatting -------------------------------------------------------------
# New cells pick up the same header/data styling already used for the rest
# of the row (bold+border header vs plain data row), overwriting the
# temporary text format applied above for J.
$ws.Range("G1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)
$ws.Range("G2").Copy()
$ws.Range("H2:N2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
